$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the company code value in A11 from ASIANPAINT to SLICE
$ws.Range("A11").Value = "SLICE"

# Update the active selection to A11 (matches final saved selection state)
$ws.Range("A11").Select()
